$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings (e.g. "6.70", "0.997")
# keep their exact textual representation instead of being parsed/rounded as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '57.916.51'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '2.571.58'
$ws.Range("E3").Value = '  -2.76%  '
$ws.Range("D5").Value = '517.13'
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("D6").Value = '142.28'
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").Value = '2.586.70'
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").Value = '0.324'
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").Value = '3.024.75'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = '57.885.75'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").Value = '20.28'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '2.524.37'
$ws.Range("E18").Value = '  -4.90%  '
$ws.Range("D19").Value = '340.13'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '4.29'
$ws.Range("E20").Value = '  -2.48%  '
$ws.Range("D21").Value = '10.22'
$ws.Range("E21").Value = '  -2.14%  '
$ws.Range("D22").Value = '6.32'
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '65.36'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '0.401'
$ws.Range("E26").Value = '  -5.47%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").Value = '2.683.09'
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").Value = '6.98'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0744'
$ws.Range("E30").Value = '  -6.49%  '
$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '6.20'
$ws.Range("E32").Value = '  -6.52%  '
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '18.68'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '149.96'
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("D36").Value = '4.00'
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("D37").Value = '1.14'
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("D38").Value = '0.870'
$ws.Range("E38").Value = '  -4.15%  '
$ws.Range("D39").Value = '35.94'
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").Value = '0.832'
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = '3.52'
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").Value = '269.82'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '10.67'
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").Value = '0.0949'
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("D47").Value = '0.586'
$ws.Range("E47").Value = '  -3.47%  '
$ws.Range("D48").Value = '18.76'
$ws.Range("E48").Value = '  -3.22%  '
$ws.Range("D49").Value = '0.0521'
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '4.65'
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.971.17'
$ws.Range("E51").Value = '  -3.30%  '

# Restore default (General) styling on the data range so no residual text-format
# style is left applied to the cells (matches original workbook formatting).
$dataRange.Style = "Normal"
